$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1527.0322
$ws.Range("I112").Value = 1050
$ws.Range("J112").Value = 1692.9565
$ws.Range("K112").Value = 3150
$ws.Range("L112").Value = 5078.8695
$ws.Range("M112").Value = -2042
$ws.Range("N112").Value = -7294.8695

$ws.Range("H127").Value = 1247.2106
$ws.Range("I127").Value = 1189.5
$ws.Range("K127").Value = 3568.5
$ws.Range("M127").Value = 1391.5

$ws.Range("H129").Value = 846.8484999999999
$ws.Range("J129").Value = 915.5909
$ws.Range("L129").Value = 2746.7727
$ws.Range("N129").Value = -12746.7727

$ws.Range("H137").Value = 3833.6667
$ws.Range("I137").Value = 1219.7037
$ws.Range("J137").Value = 7754.6113
$ws.Range("K137").Value = 3659.1111
$ws.Range("L137").Value = 23263.8339
$ws.Range("M137").Value = -1109.1111
$ws.Range("N137").Value = -28363.8339

$ws.Range("H141").Value = 5795.0835
$ws.Range("I141").Value = 3100.5144
$ws.Range("J141").Value = 100105
$ws.Range("K141").Value = 9301.5432
$ws.Range("L141").Value = 300315
$ws.Range("M141").Value = -4121.5432
$ws.Range("N141").Value = -310675


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9325.888000000001
$ws.Range("I32").Value = 8336.727999999999
$ws.Range("J32").Value = 22382.8
$ws.Range("K32").Value = 8336.727999999999
$ws.Range("L32").Value = 22382.8
$ws.Range("M32").Value = -8049.727999999999
$ws.Range("N32").Value = -22956.8

$ws.Range("H51").Value = 22247
$ws.Range("J51").Value = 22247
$ws.Range("L51").Value = 22247
$ws.Range("N51").Value = -23759

$ws.Range("H61").Value = 324266.2
$ws.Range("I61").Value = 1722.7858
$ws.Range("J61").Value = 3334671.2
$ws.Range("K61").Value = 1722.7858
$ws.Range("L61").Value = 3334671.2
$ws.Range("M61").Value = -1510.7858
$ws.Range("N61").Value = -3335095.2

$ws.Range("H74").Value = 7731.1113
$ws.Range("I74").Value = 1638.6
$ws.Range("J74").Value = 15346.75
$ws.Range("K74").Value = 1638.6
$ws.Range("L74").Value = 15346.75
$ws.Range("M74").Value = -764.5999999999999
$ws.Range("N74").Value = -17094.75

$ws.Range("H75").Value = 28000
$ws.Range("J75").Value = 28000
$ws.Range("L75").Value = 28000
$ws.Range("N75").Value = -29748

$ws.Range("H77").Value = 7731.1113
$ws.Range("I77").Value = 1638.6
$ws.Range("J77").Value = 15346.75
$ws.Range("K77").Value = 8193
$ws.Range("L77").Value = 76733.75
$ws.Range("M77").Value = -3825
$ws.Range("N77").Value = -85469.75

$ws.Range("H78").Value = 28000
$ws.Range("J78").Value = 28000
$ws.Range("L78").Value = 84000
$ws.Range("N78").Value = -92736

$ws.Range("H110").Value = 2936.5789
$ws.Range("I110").Value = 2832.9644
$ws.Range("J110").Value = 3226.7
$ws.Range("K110").Value = 2832.9644
$ws.Range("L110").Value = 3226.7
$ws.Range("M110").Value = -787.9643999999998
$ws.Range("N110").Value = -7316.7

$ws.Range("H132").Value = 5793.787
$ws.Range("I132").Value = 3941
$ws.Range("J132").Value = 11857.454
$ws.Range("K132").Value = 11823
$ws.Range("L132").Value = 35572.362
$ws.Range("M132").Value = -9293
$ws.Range("N132").Value = -40632.362

$ws.Range("H136").Value = 324266.2
$ws.Range("I136").Value = 1722.7858
$ws.Range("J136").Value = 3334671.2
$ws.Range("K136").Value = 5168.357400000001
$ws.Range("L136").Value = 10004013.6
$ws.Range("M136").Value = -2618.357400000001
$ws.Range("N136").Value = -10009113.6


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2261.87
$ws.Range("I31").Value = 1100.5217
$ws.Range("J31").Value = 3251.1667
$ws.Range("K31").Value = 1100.5217
$ws.Range("L31").Value = 3251.1667
$ws.Range("M31").Value = -805.5217
$ws.Range("N31").Value = -3841.1667

$ws.Range("H34").Value = 2261.87
$ws.Range("I34").Value = 1100.5217
$ws.Range("J34").Value = 3251.1667
$ws.Range("K34").Value = 1100.5217
$ws.Range("L34").Value = 3251.1667
$ws.Range("M34").Value = -898.5217
$ws.Range("N34").Value = -3655.1667

$ws.Range("H132").Value = 9053.809999999999
$ws.Range("I132").Value = 9006.429
$ws.Range("J132").Value = 9148.571
$ws.Range("K132").Value = 27019.287
$ws.Range("L132").Value = 27445.713
$ws.Range("M132").Value = -24489.287
$ws.Range("N132").Value = -32505.713


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 537.6667
$ws.Range("I4").Value = 308.2
$ws.Range("J4").Value = 996.6
$ws.Range("K4").Value = 924.5999999999999
$ws.Range("L4").Value = 2989.8
$ws.Range("M4").Value = -812.5999999999999
$ws.Range("N4").Value = -3213.8


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 39336.25
$ws.Range("J95").Value = 39336.25
$ws.Range("L95").Value = 39336.25
$ws.Range("N95").Value = -44828.25

$ws.Range("H132").Value = 27381.773
$ws.Range("I132").Value = 38933.566
$ws.Range("J132").Value = 2627.9285
$ws.Range("K132").Value = 116800.698
$ws.Range("L132").Value = 7883.7855
$ws.Range("M132").Value = -114270.698
$ws.Range("N132").Value = -12943.7855


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 19661
$ws.Range("I76").Value = 3500
$ws.Range("J76").Value = 25048
$ws.Range("K76").Value = 3500
$ws.Range("L76").Value = 25048
$ws.Range("M76").Value = -3162
$ws.Range("N76").Value = -25724

$ws.Range("H79").Value = 19661
$ws.Range("I79").Value = 3500
$ws.Range("J79").Value = 25048
$ws.Range("K79").Value = 3500
$ws.Range("L79").Value = 25048
$ws.Range("M79").Value = -2330
$ws.Range("N79").Value = -27388

$ws.Range("H132").Value = 30392.8
$ws.Range("I132").Value = 12718.546
$ws.Range("K132").Value = 38155.638
$ws.Range("M132").Value = -35625.638


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4039.625
$ws.Range("I132").Value = 4674.657
$ws.Range("J132").Value = 2329.923
$ws.Range("K132").Value = 14023.971
$ws.Range("L132").Value = 6989.768999999999
$ws.Range("M132").Value = -11493.971
$ws.Range("N132").Value = -12049.769

$ws.Range("H136").Value = 10168.083
$ws.Range("I136").Value = 18867.834
$ws.Range("J136").Value = 1468.3334
$ws.Range("K136").Value = 56603.50199999999
$ws.Range("L136").Value = 4405.0002
$ws.Range("M136").Value = -54053.50199999999
$ws.Range("N136").Value = -9505.0002

